# Update the "F" column numeric values (view counts) on the "展览" and
# "全部类型" worksheets to reflect newly generated data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 3622
$ws1.Range("F16").Value = 630
$ws1.Range("F17").Value = 127
$ws1.Range("F18").Value = 811
$ws1.Range("F19").Value = 30
$ws1.Range("F25").Value = 2919
$ws1.Range("F26").Value = 5348
$ws1.Range("F32").Value = 313
$ws1.Range("F37").Value = 150
$ws1.Range("F45").Value = 457

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 3622
$ws4.Range("F17").Value = 630
$ws4.Range("F18").Value = 127
$ws4.Range("F19").Value = 811
$ws4.Range("F20").Value = 30
$ws4.Range("F26").Value = 2919
$ws4.Range("F27").Value = 5348
$ws4.Range("F33").Value = 313
$ws4.Range("F38").Value = 150
$ws4.Range("F46").Value = 457
